$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string header text updates ---
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# --- Simple same-format numeric value updates (crime stat table) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 13
$ws.Range("H16").Value = 18.181818181818
$ws.Range("I16").Value = 203
$ws.Range("J16").Value = 106
$ws.Range("K16").Value = 91.509433962264
$ws.Range("L16").Value = 48.175182481751
$ws.Range("M16").Value = 41.958041958042
$ws.Range("N16").Value = -76.746849942726
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 5
$ws.Range("I17").Value = 247
$ws.Range("J17").Value = 230
$ws.Range("K17").Value = 7.391304347826
$ws.Range("L17").Value = 15.962441314554
$ws.Range("M17").Value = 133.018867924528
$ws.Range("N17").Value = -19.281045751634
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 112.5
$ws.Range("I18").Value = 237
$ws.Range("J18").Value = 147
$ws.Range("K18").Value = 61.224489795918
$ws.Range("L18").Value = 16.176470588235
$ws.Range("M18").Value = -4.048582995951
$ws.Range("N18").Value = -88.461538461538
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = -45.945945945945
$ws.Range("I19").Value = 715
$ws.Range("J19").Value = 460
$ws.Range("K19").Value = 55.434782608695
$ws.Range("L19").Value = 55.434782608695
$ws.Range("M19").Value = 82.864450127877
$ws.Range("N19").Value = -49.326718639262
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 14
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 180
$ws.Range("J20").Value = 90
$ws.Range("K20").Value = 100
$ws.Range("L20").Value = 28.571428571428
$ws.Range("M20").Value = -1.098901098901
$ws.Range("N20").Value = -94.352055224348
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -24
$ws.Range("F21").Value = 105
$ws.Range("H21").Value = -13.223140495867
$ws.Range("I21").Value = 1593
$ws.Range("J21").Value = 1042
$ws.Range("K21").Value = 52.879078694817
$ws.Range("L21").Value = 35.921501706484
$ws.Range("M21").Value = 47.636700648748
$ws.Range("N21").Value = -79.737980157720
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = -43.75
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 6
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 99
$ws.Range("J23").Value = 65
$ws.Range("K23").Value = 52.307692307692
$ws.Range("L23").Value = 73.684210526315
$ws.Range("M23").Value = 241.379310344828
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -4.761904761904
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -1.960784313725
$ws.Range("I24").Value = 1555
$ws.Range("J24").Value = 1049
$ws.Range("K24").Value = 48.236415633937
$ws.Range("L24").Value = 42.791551882461
$ws.Range("M24").Value = 46.284101599247
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = -9.090909090909
$ws.Range("I25").Value = 411
$ws.Range("J25").Value = 340
$ws.Range("K25").Value = 20.882352941176
$ws.Range("L25").Value = 31.309904153354
$ws.Range("M25").Value = -16.632860040568
$ws.Range("G26").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 48
$ws.Range("J27").Value = 47
$ws.Range("K27").Value = 2.127659574468
$ws.Range("L27").Value = 84.615384615384
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("G30").Value = 2

# --- Cells switching FROM the text "N/A" placeholder TO a real number ---
# (set NumberFormat first so the new value is stored as a genuine number
#  with the right display style, matching the neighboring numeric cells)
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = -100
$ws.Range("G22").NumberFormat = '#,##0'
$ws.Range("G22").Value = 1
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H22").Value = -100
$ws.Range("C23").NumberFormat = '#,##0'
$ws.Range("C23").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = -50

# --- Cells switching FROM a real number TO the text "N/A" placeholder ---
# (force Text format so the numeric-looking literal is stored as a string,
#  then restore General format so the cell style matches its text neighbors)
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C28").NumberFormat = "general"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C29").NumberFormat = "general"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D30").NumberFormat = "general"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E30").NumberFormat = "general"
